# Add new columns I ("I0") and J ("IF") per commit "I0 and IF added".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns: copy formatting from the existing
# header cell (H1) so the new headers match the bold/centered/bordered style,
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..81: new numeric values for columns I and J.
$data = @(
    @(2, 3, 4),
    @(3, 5, 6),
    @(4, 5, 6),
    @(5, 8, 8),
    @(6, 6, 6),
    @(7, 6, 6),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 7, 8),
    @(11, 8, 9),
    @(12, 10, 10),
    @(13, 7, 7),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 5, 6),
    @(17, 11, 11),
    @(18, 8, 8),
    @(19, 6, 6),
    @(20, 6, 6),
    @(21, 10, 10),
    @(22, 7, 8),
    @(23, 5, 6),
    @(24, 8, 8),
    @(25, 6, 7),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 6, 6),
    @(29, 7, 7),
    @(30, 7, 7),
    @(31, 7, 7),
    @(32, 9, 9),
    @(33, 6, 6),
    @(34, 5, 5),
    @(35, 6, 6),
    @(36, 7, 7),
    @(37, 6, 6),
    @(38, 6, 6),
    @(39, 6, 7),
    @(40, 9, 9),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 7, 7),
    @(44, 4, 4),
    @(45, 7, 7),
    @(46, 7, 7),
    @(47, 12, 12),
    @(48, 6, 7),
    @(49, 6, 7),
    @(50, 9, 9),
    @(51, 6, 6),
    @(52, 6, 6),
    @(53, 6, 6),
    @(54, 6, 7),
    @(55, 9, 9),
    @(56, 6, 7),
    @(57, 7, 7),
    @(58, 8, 8),
    @(59, 7, 7),
    @(60, 9, 9),
    @(61, 7, 8),
    @(62, 5, 6),
    @(63, 6, 6),
    @(64, 7, 8),
    @(65, 4, 5),
    @(66, 4, 5),
    @(67, 8, 8),
    @(68, 6, 6),
    @(69, 5, 5),
    @(70, 6, 6),
    @(71, 5, 5),
    @(72, 6, 6),
    @(73, 5, 5),
    @(74, 8, 8),
    @(75, 4, 5),
    @(76, 3, 3),
    @(77, 8, 8),
    @(78, 9, 9),
    @(79, 8, 8),
    @(80, 7, 7),
    @(81, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
